$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the current extent of the data (header row + data rows,
# columns A..S i.e. up to the "04-02_A"/"04-02_0" pair).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# The last two columns (R, S) hold the most recent date pair
# ("04-02_A" / "04-02_0"). We duplicate that pair into two new
# columns (T, U) to represent the next day's pair ("04-03_A" / "04-03_0").
$srcAcol = $lastCol - 1
$srcOcol = $lastCol
$dstAcol = $lastCol + 1
$dstOcol = $lastCol + 2

$srcRange = $ws.Range($ws.Cells.Item(1, $srcAcol), $ws.Cells.Item($lastRow, $srcOcol))
$dstRange = $ws.Cells.Item(1, $dstAcol)

# Copy formatting + values of the R:S block into the new T:U block.
$srcRange.Copy($dstRange)

# Fix up the new header labels (copy duplicated the old "04-02_*" text).
$ws.Cells.Item(1, $dstAcol).Value2 = "04-03_A"
$ws.Cells.Item(1, $dstOcol).Value2 = "04-03_0"

# The "_0" column (S) used to be stored as text (inline string) even when
# it held a numeric attendance total. Convert any non-blank S cell back
# into a real number, matching the rest of the "_0" columns.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $srcOcol)
    $val = $cell.Value2
    if ($val -ne $null -and "$val" -ne "") {
        $cell.Value2 = $val + 0
    }
}
